$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell, far away from the used range (B2:C39), used to stage new
# text so that writing the literal string value does not disturb the
# destination cell's existing format (quotePrefix / wrap, etc. - this
# runtime's Value setter resets those on the written cell).
$scratch = $ws.Range("Z50")

function Set-TextPreserveStyle($targetAddress, $text) {
    $target = $ws.Range($targetAddress)
    $scratch.Value = $text
    $scratch.Copy()
    $target.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
    $target.Copy()
    $target.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    $scratch.ClearContents()
    # Undo any row-height autofit the scratch write left behind.
    $ws.Rows.Item(50).AutoFit()
}

# --- Day 1 block becomes 13/10 (was 09/10) ---

# New date
$ws.Range("B4").Value = 41925

# Ke hoach (plan) -> new shared string, added first
Set-TextPreserveStyle "C5" "- Công việc 1: Chỉnh sửa để code đã có chạy.
- Công việc 2: Biuld màn hình hiển thị cho phần chức năng `" Công việc `" ( thêm công việc, tất cả công việc, công việc đã giao, công việc được giao)."

# Trang thai (status)
$ws.Range("C7").Value = "Chậm tiến độ"

# Van de gap phai (problem encountered)
Set-TextPreserveStyle "C8" "- Vấn đề 1: option menu ở FragmentCongViec khi kéo NavigationDrawer ra vẫn chưa ẩn đi được."

# Giai quyet van de (solution)
Set-TextPreserveStyle "C10" "- Vấn đề 1: đang tìm hiểu vẫn chưa có phương hướng cụ thể."

# Ket qua dat duoc (results achieved)
Set-TextPreserveStyle "C6" "- Công việc 1: hoàn thành.
- Công việc 2: 50%"

# Ke hoach ngay mai (plan for tomorrow)
$ws.Range("C12").Value = "Tiếp tục công việc 2 , giải quyết vấn đề 1 và code để đổ data ra listview"

$ws.Range("C12").Select()

$wb.Save()
